$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.41"
$ws.Range("E2").Value = "'-0.05%"
$ws.Range("D3").Value = "'41.09"
$ws.Range("E3").Value = "'-0.16%"
$ws.Range("D4").Value = "'5.178"
$ws.Range("E4").Value = "'0.89%"
$ws.Range("D5").Value = "'0.07665"
$ws.Range("E5").Value = "'0.48%"
$ws.Range("D6").Value = "'1.710"
$ws.Range("E6").Value = "'5.62%"
$ws.Range("D7").Value = "'0.9150"
$ws.Range("E8").Value = "'-2.09%"
$ws.Range("E9").Value = "'11.09%"
$ws.Range("D10").Value = "'0.1815"
$ws.Range("E10").Value = "'1.39%"
$ws.Range("D11").Value = "'0.09161"
$ws.Range("E11").Value = "'0.92%"
$ws.Range("D12").Value = "'0.04193"
$ws.Range("E12").Value = "'-1.44%"
$ws.Range("E13").Value = "'0.21%"
$ws.Range("D14").Value = "'0.001304"
$ws.Range("E14").Value = "'4.14%"
$ws.Range("D15").Value = "'0.005750"
$ws.Range("E15").Value = "'2.24%"
$ws.Range("D16").Value = "'3.344"
$ws.Range("E16").Value = "'0.04%"
$ws.Range("E17").Value = "'1.20%"
$ws.Range("D19").Value = "'7.401"
$ws.Range("E19").Value = "'11.03%"
$ws.Range("D20").Value = "'0.1357"
$ws.Range("E20").Value = "'-0.51%"
$ws.Range("D21").Value = "'0.2729"
$ws.Range("E21").Value = "'0.81%"
$ws.Range("D22").Value = "'0.04024"
$ws.Range("E22").Value = "'-0.06%"
$ws.Range("E23").Value = "'0.60%"
$ws.Range("D24").Value = "'0.004076"
$ws.Range("E24").Value = "'0.32%"
$ws.Range("E25").Value = "'0.31%"
$ws.Range("D38").Value = "'0.02526"
$ws.Range("E38").Value = "'3.90%"
$ws.Range("D39").Value = "'0.05296"
$ws.Range("E39").Value = "'1.05%"
$ws.Range("D40").Value = "'0.007844"
$ws.Range("E40").Value = "'0.60%"
$ws.Range("D41").Value = "'0.1309"
$ws.Range("E41").Value = "'0.55%"
$ws.Range("D42").Value = "'0.006661"
$ws.Range("E42").Value = "'-5.48%"
$ws.Range("D43").Value = "'0.001875"
$ws.Range("E43").Value = "'-3.80%"
$ws.Range("D44").Value = "'0.008139"
$ws.Range("E44").Value = "'-3.54%"
$ws.Range("D45").Value = "'0.3070"
$ws.Range("E45").Value = "'-8.25%"
$ws.Range("D46").Value = "'0.00006807"
$ws.Range("E46").Value = "'2.57%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.44%"
$ws.Range("D48").Value = "'0.2244"
$ws.Range("E48").Value = "'308.76%"
$ws.Range("D50").Value = "'0.00002108"
$ws.Range("E50").Value = "'0.44%"
$ws.Range("D51").Value = "'0.0002008"
$ws.Range("E51").Value = "'0.44%"
